$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

# The "users" sheet had accidentally leaked everyone's bcrypt password hash
# ("$2a$04$AvIPbbagfXqklLVaXXYiZ.xHYIDmL.xvvd.UrHySr5k4AEaaNH/82") in column G
# (rows 2-101). Scrub it by overwriting those cells with the generic
# "password" placeholder text (the same text already used for the column
# header in G1), so the real hash no longer appears anywhere in the workbook.
$ws.Range("G2:G101").Value = "password"

# Restore the sheet's scroll/selection state as left by the author.
$ws.Activate()
$ws.Range("G37").Select()
